$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("#7cc867#f9cd59#c885da#fb5b89", $true, $false, $false, $false, $false, $true, 1, $false, "#7cc867: 12`r#f9cd59: 16`r#c885da: 12`r#fb5b89: 12", 2)
